$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet held one monthly M2 data point per row, but was missing the
# Aug-2022 and Sep-2022 rows, and had a stale Jul-2022 value. Shift the
# existing Oct-2022..Jan-2023 rows (369:372) down to 371:374 first -- using
# Copy so the date/number formatting (style) carried by those rows moves
# with them -- then fill in the corrected Jul-2022 row and the two newly
# inserted rows.

$ws.Range("A369:G372").Copy($ws.Range("A371"))

# Row 368: Jul-2022 open/high/low/close corrected.
$ws.Range("C368:F368").Value = 192778000000

# Row 369 (new): Aug-2022 data point.
$ws.Range("A369").Value = 44774.41666666666
$ws.Range("B369").Value = "ECONOMICS:PEM2"
$ws.Range("C369:F369").Value = 295831000000
$ws.Range("G369").Value = 0

# Row 370 (new): Sep-2022 data point.
$ws.Range("A370").Value = 44805.41666666666
$ws.Range("B370").Value = "ECONOMICS:PEM2"
$ws.Range("C370:F370").Value = 295831000000
$ws.Range("G370").Value = 0

$ws.Range("A1").Select()
